$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Cylinder" table (rows 5-8) needs to be sorted by column A (cyl count)
# ascending. Before: 6, 4, 6, 8. After: 4, 6, 6, 8.
# That only requires swapping rows 5 and 6; row 7 and row 8 already sit in
# their final, sorted position.

# --- capture the values of rows 5 and 6 (columns A-G) before overwriting ---
$A5 = $ws.Range("A5").Value(); $B5 = $ws.Range("B5").Value(); $C5 = $ws.Range("C5").Value()
$D5 = $ws.Range("D5").Value(); $E5 = $ws.Range("E5").Value(); $F5 = $ws.Range("F5").Value(); $G5 = $ws.Range("G5").Value()

$A6 = $ws.Range("A6").Value(); $B6 = $ws.Range("B6").Value(); $C6 = $ws.Range("C6").Value()
$D6 = $ws.Range("D6").Value(); $E6 = $ws.Range("E6").Value(); $F6 = $ws.Range("F6").Value(); $G6 = $ws.Range("G6").Value()

# --- row 5 becomes what used to be row 6 ---
$ws.Range("A5").Value = $A6
$ws.Range("B5").Value = $B6
$ws.Range("C5").Value = $C6
$ws.Range("D5").Value = $D6
$ws.Range("E5").Value = $E6
$ws.Range("F5").Value = $F6
$ws.Range("G5").Value = $G6

# --- row 6 becomes what used to be row 5 ---
$ws.Range("A6").Value = $A5
$ws.Range("B6").Value = $B5
$ws.Range("C6").Value = $C5
$ws.Range("D6").Value = $D5
$ws.Range("E6").Value = $E5
$ws.Range("F6").Value = $F5
$ws.Range("G6").Value = $G5

# Rows 6 and 7 now share the same cylinder count (6), so column A is merged
# across the two rows, with the value top-aligned in the merged cell.
$ws.Range("A6:A7").Merge()
$ws.Range("A6").VerticalAlignment = -4160   # xlTop

# The (now redundant) value in A7 is cleared, and its formatting is reset to
# match the plain, unstyled look used by the other "spacer" cells (e.g. H7).
$ws.Range("A7").ClearContents()
$ws.Range("H7").Copy()
$ws.Range("A7").PasteSpecial(-4122)         # xlPasteFormats
